$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value2 = '57.340.93'
$ws.Range('E2').Value2 = '  +5.02%  '

# Row 3
$ws.Range('D3').Value2 = '2.360.87'
$ws.Range('E3').Value2 = '  +2.97%  '

# Row 4
$ws.Range('E4').Value2 = '  -0.18%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value2 = '520.03'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value2 = '  +3.25%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value2 = '134.59'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value2 = '  +3.28%  '

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value2 = '1.00'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value2 = '  +0.25%  '

# Row 8
$ws.Range('E8').Value2 = '  +2.06%  '

# Row 9
$ws.Range('D9').Value2 = '2.357.95'
$ws.Range('E9').Value2 = '  +1.75%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value2 = '0.103'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value2 = '  +7.42%  '

# Row 11
$ws.Range('E11').Value2 = '  +1.09%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value2 = '5.22'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value2 = '  +6.57%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value2 = '0.343'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value2 = '  +1.82%  '

# Row 14
$ws.Range('B14').Value2 = 'WrappedliquidstakedEther2.0'
$ws.Range('C14').Value2 = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D14').Value2 = '2.780.87'
$ws.Range('E14').Value2 = '  +2.89%  '

# Row 15
$ws.Range('B15').Value2 = 'Avalanche'
$ws.Range('C15').Value2 = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value2 = '23.87'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value2 = '  +2.68%  '

# Row 16
$ws.Range('D16').Value2 = '57.184.98'
$ws.Range('E16').Value2 = '  +4.41%  '

# Row 17
$ws.Range('E17').Value2 = '  +3.63%  '

# Row 18
$ws.Range('D18').Value2 = '2.351.41'
$ws.Range('E18').Value2 = '  +1.90%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value2 = '10.56'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value2 = '  +2.04%  '

# Row 20
$ws.Range('E20').Value2 = '  +3.06%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value2 = '321.44'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value2 = '  +4.72%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value2 = '6.72'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value2 = '  +5.80%  '

# Row 23
$ws.Range('E23').Value2 = '  -0.04%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value2 = '61.40'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value2 = '  +0.99%  '

# Row 25
$ws.Range('E25').Value2 = '  +7.32%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value2 = '0.996'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value2 = '  +0.33%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value2 = '7.81'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value2 = '  +5.34%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value2 = '172.53'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value2 = '  -0.99%  '

# Row 29
$ws.Range('E29').Value2 = '  +9.31%  '

# Row 30
$ws.Range('D30').Value2 = '0.0₃0741'
$ws.Range('E30').Value2 = '  +3.21%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value2 = '6.30'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value2 = '  +4.13%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value2 = '1.68'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value2 = '  +3.39%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value2 = '18.41'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value2 = '  +2.49%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value2 = '1.00'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value2 = '  +0.10%  '

# Row 35
$ws.Range('E35').Value2 = '  +2.16%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value2 = '1.00'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value2 = '  +0.48%  '

# Row 37
$ws.Range('E37').Value2 = '  +4.32%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value2 = '4.04'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value2 = '  +7.28%  '

# Row 39
$ws.Range('B39').Value2 = 'OKB'
$ws.Range('C39').Value2 = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value2 = '37.66'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value2 = '  +4.06%  '

# Row 40
$ws.Range('B40').Value2 = 'Stacks'
$ws.Range('C40').Value2 = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value2 = '1.52'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value2 = '  +6.92%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value2 = '0.384'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value2 = '  +1.69%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value2 = '139.60'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value2 = '  +11.13%  '

# Row 43
$ws.Range('E43').Value2 = '  +5.62%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value2 = '279.11'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value2 = '  +11.96%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value2 = '5.17'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value2 = '  +2.99%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value2 = '0.0511'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value2 = '  +3.11%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value2 = '0.0931'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value2 = '  +3.60%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value2 = '0.564'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value2 = '  +1.87%  '

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value2 = '0.383'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value2 = '  +1.85%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value2 = '0.0216'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value2 = '  +4.21%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value2 = '17.00'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value2 = '  +2.41%  '

